$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.206015333333333
$ws.Range("H2").Value = 18.618046
$ws.Range("I2").Value = 0.0150172404156507
$ws.Range("J2").Value = 0.0150172404156507
$ws.Range("M2").Value = 0.2303363333333333
$ws.Range("N2").Value = 0.691009
$ws.Range("O2").Value = 0.0420565315194687
$ws.Range("P2").Value = 0.0420565315194687
$ws.Range("Q2").Value = 1.429470816490444
$ws.Range("R2").Value = 12.865237348414
$ws.Range("S2").Value = 0.0006315730448762529
$ws.Range("T2").Value = 0.0006315730448762527
$ws.Range("G3").Value = 6.206015333333333
$ws.Range("H3").Value = 18.618046
$ws.Range("I3").Value = 0.0150172404156507
$ws.Range("J3").Value = 0.0150172404156507
$ws.Range("O3").Value = 0.8440851393264226
$ws.Range("P3").Value = 0.8440851393264227
$ws.Range("Q3").Value = 28.68983793258933
$ws.Range("R3").Value = 258.208541393304
$ws.Range("S3").Value = 0.01267582946854291
$ws.Range("T3").Value = 0.01267582946854291
$ws.Range("G4").Value = 6.206015333333333
$ws.Range("H4").Value = 18.618046
$ws.Range("I4").Value = 0.0150172404156507
$ws.Range("J4").Value = 0.0150172404156507
$ws.Range("M4").Value = 0.6235823333333333
$ws.Range("N4").Value = 1.870747
$ws.Range("O4").Value = 0.1138583291541087
$ws.Range("P4").Value = 0.1138583291541087
$ws.Range("Q4").Value = 3.869961522262444
$ws.Range("R4").Value = 34.829653700362
$ws.Range("S4").Value = 0.001709837902231542
$ws.Range("T4").Value = 0.001709837902231542
$ws.Range("I5").Value = 0.9317452840597572
$ws.Range("J5").Value = 0.9317452840597571
$ws.Range("M5").Value = 0.2303363333333333
$ws.Range("N5").Value = 0.691009
$ws.Range("O5").Value = 0.0420565315194687
$ws.Range("P5").Value = 0.0420565315194687
$ws.Range("Q5").Value = 88.69157415752211
$ws.Range("R5").Value = 798.224167417699
$ws.Range("S5").Value = 0.03918597490717549
$ws.Range("T5").Value = 0.03918597490717549
$ws.Range("I6").Value = 0.9317452840597572
$ws.Range("J6").Value = 0.9317452840597571
$ws.Range("O6").Value = 0.8440851393264226
$ws.Range("P6").Value = 0.8440851393264227
$ws.Range("S6").Value = 0.7864723479123173
$ws.Range("T6").Value = 0.7864723479123173
$ws.Range("I7").Value = 0.9317452840597572
$ws.Range("J7").Value = 0.9317452840597571
$ws.Range("M7").Value = 0.6235823333333333
$ws.Range("N7").Value = 1.870747
$ws.Range("O7").Value = 0.1138583291541087
$ws.Range("P7").Value = 0.1138583291541087
$ws.Range("Q7").Value = 240.1119179062241
$ws.Range("R7").Value = 2161.007261156018
$ws.Range("S7").Value = 0.1060869612402644
$ws.Range("T7").Value = 0.1060869612402644
$ws.Range("G8").Value = 22.00088566666667
$ws.Range("H8").Value = 66.002657
$ws.Range("I8").Value = 0.05323747552459213
$ws.Range("J8").Value = 0.05323747552459213
$ws.Range("M8").Value = 0.2303363333333333
$ws.Range("N8").Value = 0.691009
$ws.Range("O8").Value = 0.0420565315194687
$ws.Range("P8").Value = 0.0420565315194687
$ws.Range("Q8").Value = 5.067603334545889
$ws.Range("R8").Value = 45.608430010913
$ws.Range("S8").Value = 0.002238983567416952
$ws.Range("T8").Value = 0.002238983567416952
$ws.Range("G9").Value = 22.00088566666667
$ws.Range("H9").Value = 66.002657
$ws.Range("I9").Value = 0.05323747552459213
$ws.Range("J9").Value = 0.05323747552459213
$ws.Range("O9").Value = 0.8440851393264226
$ws.Range("P9").Value = 0.8440851393264227
$ws.Range("Q9").Value = 101.7080703555187
$ws.Range("R9").Value = 915.372633199668
$ws.Range("S9").Value = 0.04493696194556236
$ws.Range("T9").Value = 0.04493696194556236
$ws.Range("G10").Value = 22.00088566666667
$ws.Range("H10").Value = 66.002657
$ws.Range("I10").Value = 0.05323747552459213
$ws.Range("J10").Value = 0.05323747552459213
$ws.Range("M10").Value = 0.6235823333333333
$ws.Range("N10").Value = 1.870747
$ws.Range("O10").Value = 0.1138583291541087
$ws.Range("P10").Value = 0.1138583291541087
$ws.Range("Q10").Value = 13.71936361941989
$ws.Range("R10").Value = 123.474272574779
$ws.Range("S10").Value = 0.006061530011612817
$ws.Range("T10").Value = 0.006061530011612817
